$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.556.71'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '2.162.98'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.85'
$ws.Range("E5").Value = '  -1.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.621'
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '62.59'
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("E9").Value = '  -1.03%  '
$ws.Range("E10").Value = '  -0.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.103'
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("E12").Value = '  -1.65%  '
$ws.Range("D13").Value = '2.484.12'
$ws.Range("E13").Value = '  -0.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.67'
$ws.Range("E14").Value = '  -2.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.805'
$ws.Range("E15").Value = '  -1.99%  '
$ws.Range("E16").Value = '  -1.65%  '
$ws.Range("D17").Value = '2.176.84'
$ws.Range("E17").Value = '  +0.36%  '
$ws.Range("D18").Value = '39.555.60'
$ws.Range("E18").Value = '  +0.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.54'
$ws.Range("E19").Value = '  -0.82%  '
$ws.Range("D20").Value = '0.0₃0885'
$ws.Range("E20").Value = '  +3.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.99'
$ws.Range("E21").Value = '  -2.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.50'
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.35'
$ws.Range("E24").Value = '  +0.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.30'
$ws.Range("E25").Value = '  -4.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.38'
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.43'
$ws.Range("E27").Value = '  -2.98%  '
$ws.Range("E28").Value = '  +0.19%  '
$ws.Range("E29").Value = '  +1.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.61'
$ws.Range("E30").Value = '  -0.47%  '
$ws.Range("E31").Value = '  +4.46%  '
$ws.Range("E32").Value = '  +0.48%  '
$ws.Range("E33").Value = '  -3.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.68'
$ws.Range("E34").Value = '  -3.15%  '
$ws.Range("E35").Value = '  -3.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0616'
$ws.Range("E36").Value = '  -0.54%  '
$ws.Range("E37").Value = '  +7.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.38'
$ws.Range("E38").Value = '  -2.15%  '
$ws.Range("E39").Value = '  +0.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.89'
$ws.Range("E40").Value = '  +17.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '101.98'
$ws.Range("E41").Value = '  -1.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0226'
$ws.Range("E42").Value = '  -1.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.67'
$ws.Range("E43").Value = '  -3.05%  '
$ws.Range("D44").Value = '1.512.01'
$ws.Range("E44").Value = '  -1.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.20'
$ws.Range("E45").Value = '  +0.66%  '
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0914'
$ws.Range("E48").Value = '  -0.92%  '
$ws.Range("E49").Value = '  -1.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000197'
$ws.Range("E50").Value = '  +33.25%  '
$ws.Range("E51").Value = '  -0.03%  '
